# /tmp/work/edit.ps1
# Applies the "Actualiza base de datos EC y agrega parte 1 de nuevos estado
# de cuenta" commit: bumps the overdue-balance figure, the period count,
# re-sorts the overdue-period detail table into ascending order, and
# appends one more overdue period ("2508") at the bottom - pushing the
# trailing signature block down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# ------------------------------------------------------------------
# 1) Header figures: "VALOR MORA" and "Cant. Periodos" both increase
#    now that a new overdue period has been added to the table below.
# ------------------------------------------------------------------
$ws.Range("E11").Value = 4620000
$ws.Range("F13").Value = 110

# ------------------------------------------------------------------
# 2) Make room for the new period row. The detail table currently runs
#    from row 16 down to row 124 (the last row, which carries a
#    heavier closing border). Inserting a row at 125 pushes the blank
#    spacer rows and the closing signature block (rows 129-130) down
#    to rows 130-131, matching the new layout.
# ------------------------------------------------------------------
$ws.Rows("125:125").Insert()

# The insert leaves row 124 still wearing the special "closing" border
# style and row 125 with a blank/default style. Move that closing style
# onto the new last row (125), and restore row 124 to the regular
# interior-row style (copied from row 123).
$ws.Range("B124:J124").Copy() | Out-Null
$ws.Range("B125:J125").PasteSpecial(-4122) | Out-Null

$ws.Range("B123:J123").Copy() | Out-Null
$ws.Range("B124:J124").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 3) Re-write the "Periodo Mora" column for the whole table (rows 16-125)
#    in ascending chronological order, from 1607 through the newly
#    added 2508. Worker identity and debt figures (columns B,C,D,F,G)
#    are identical on every row already and are left untouched.
# ------------------------------------------------------------------
$periods = @(
    "1607","1608","1609","1610","1611","1612",
    "1701","1702","1703","1704","1705","1706","1707","1708","1709","1710","1711","1712",
    "1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812",
    "1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912",
    "2001","2002","2003","2004","2005","2006","2007","2008","2009","2010","2011","2012",
    "2101","2102","2103","2104","2105","2106","2107","2108","2109","2110","2111","2112",
    "2201","2202","2203","2204","2205","2206","2207","2208","2209","2210","2211","2212",
    "2301","2302","2303","2304","2305","2306","2307","2308","2309","2310","2311","2312",
    "2401","2402","2403","2404","2405","2406","2407","2408","2409","2410","2411","2412",
    "2501","2502","2503","2504","2505","2506","2507","2508"
)

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 2).Value = "CC"
    $ws.Cells.Item($row, 3).Value = "73558563"
    $ws.Cells.Item($row, 4).Value = "WALFREDO ALVEAR MARRUGO"
    $ws.Cells.Item($row, 5).Value = $periods[$i]
    $ws.Cells.Item($row, 6).Value = 42000
    $ws.Cells.Item($row, 7).Value = 1050000
}
